$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    "Saquon Barkley",
    "Kyler Murray",
    "Davante Adams",
    "Dak Prescott",
    "Jonathan Taylor",
    "Alvin Kamara",
    "Tyreek Hill",
    "Lamar Jackson",
    "DK Metcalf",
    "A.J. Brown",
    "Derrick Henry",
    "Justin Jefferson",
    "Nick Chubb",
    "Stefon Diggs",
    "Travis Kelce"
)

$row = 2
foreach ($player in $players) {
    $ws.Cells.Item($row, 1).Value = $player
    $row = $row + 1
}

$ws.Range("A2:A16").Select()
